$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16512555862088292"
$ws1.Range("B2").Value = "go_stims-16512555861668358.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555861918302.csv"
$ws1.Range("B4").Value = "go_stims-1651255586193826.csv"
$ws1.Range("B5").Value = "GNG_stims-1651255586207827.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16512555888672953"
$ws2.Range("B2").Value = "ZB-match_4-16512555866974907.csv"
$ws2.Range("B3").Value = "TB-16512555874864936.csv"
$ws2.Range("B4").Value = "OB-16512555873224885.csv"
$ws2.Range("B5").Value = "OB-16512555868754902.csv"
$ws2.Range("B6").Value = "ZB-match_4-1651255586235833.csv"
$ws2.Range("B7").Value = "TB-1651255588852296.csv"
$ws2.Range("B8").Value = "TB-16512555876914897.csv"
$ws2.Range("B9").Value = "ZB-match_1-1651255586572489.csv"
$ws2.Range("B10").Value = "OB-16512555871144881.csv"

# --- Sheet 3: RS_TO (name only) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16512555888743057"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1651255588931411"
$ws4.Range("B2").Value = "MM_stims-16512555888995714.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555888782966.csv"
$ws4.Range("B4").Value = "MM_stims-16512555889151974.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555888995714.csv"
$ws4.Range("B6").Value = "MM_stims-16512555889304163.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555889151974.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1651255589007733"
$ws5.Range("B2").Value = "vSAT_stims-1651255588969415.csv"
$ws5.Range("B3").Value = "SAT_stims-16512555889532733.csv"
$ws5.Range("B4").Value = "SAT_stims-16512555889408152.csv"
$ws5.Range("B5").Value = "vSAT_stims-1651255588992604.csv"
